$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.693.99'
$ws.Range("E2").Value = '  +1.60%  '

$ws.Range("D3").Value = '3.187.36'
$ws.Range("E3").Value = '  -0.54%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").Value = '''592.90'
$ws.Range("E5").Value = '  +0.08%  '

$ws.Range("D6").Value = '''137.05'
$ws.Range("E6").Value = '  +0.57%  '

$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").Value = '3.186.42'
$ws.Range("E8").Value = '  -0.60%  '

$ws.Range("D9").Value = '''0.515'
$ws.Range("E9").Value = '  +1.53%  '

$ws.Range("E10").Value = '  -0.72%  '

$ws.Range("D11").Value = '''5.40'
$ws.Range("E11").Value = '  +1.15%  '

$ws.Range("E12").Value = '  +0.64%  '

$ws.Range("E13").Value = '  +1.44%  '

$ws.Range("D14").Value = '''35.05'
$ws.Range("E14").Value = '  +4.47%  '

$ws.Range("D15").Value = '3.709.62'
$ws.Range("E15").Value = '  -0.71%  '

$ws.Range("E16").Value = '  -0.55%  '

$ws.Range("D17").Value = '3.185.71'
$ws.Range("E17").Value = '  -0.58%  '

$ws.Range("D18").Value = '63.655.94'
$ws.Range("E18").Value = '  +1.34%  '

$ws.Range("D19").Value = '''6.59'
$ws.Range("E19").Value = '  -1.63%  '

$ws.Range("D20").Value = '''463.52'

$ws.Range("D21").Value = '''14.01'
$ws.Range("E21").Value = '  +0.92%  '

$ws.Range("D22").Value = '''0.702'
$ws.Range("E22").Value = '  -1.66%  '

$ws.Range("D23").Value = '''7.70'
$ws.Range("E23").Value = '  +0.28%  '

$ws.Range("B24").Value = 'InternetComputer(DFINITY)'
$ws.Range("C24").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D24").Value = '''13.28'
$ws.Range("E24").Value = '  -0.71%  '

$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = '''83.46'
$ws.Range("E25").Value = '  -0.99%  '

$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("E27").Value = '  +0.16%  '

$ws.Range("D29").Value = '''6.93'
$ws.Range("E29").Value = '  +0.15%  '

$ws.Range("E30").Value = '  +0.70%  '

$ws.Range("D31").Value = '''7.81'
$ws.Range("E31").Value = '  -0.93%  '

$ws.Range("D32").Value = '''27.57'
$ws.Range("E32").Value = '  -0.16%  '

$ws.Range("E33").Value = '  -0.34%  '

$ws.Range("E34").Value = '  +0.88%  '

$ws.Range("E35").Value = '  -1.51%  '

$ws.Range("D36").Value = '''5.92'
$ws.Range("E36").Value = '  +1.10%  '

$ws.Range("D37").Value = '0.0₃0740'
$ws.Range("E37").Value = '  +5.93%  '

$ws.Range("D38").Value = '''51.71'
$ws.Range("E38").Value = '  -0.08%  '

$ws.Range("E39").Value = '  +0.64%  '

$ws.Range("D40").Value = '''8.19'
$ws.Range("E40").Value = '  +1.03%  '

$ws.Range("E41").Value = '  -0.80%  '

$ws.Range("D42").Value = '''2.68'
$ws.Range("E42").Value = '  +2.15%  '

$ws.Range("D43").Value = '''398.04'
$ws.Range("E43").Value = '  -5.13%  '

$ws.Range("D44").Value = '2.794.49'
$ws.Range("E44").Value = '  -7.27%  '

$ws.Range("E45").Value = '  -0.23%  '

$ws.Range("E46").Value = '  +0.44%  '

$ws.Range("D47").Value = '''127.97'
$ws.Range("E47").Value = '  +2.31%  '

$ws.Range("D49").Value = '''35.76'
$ws.Range("E49").Value = '  +0.15%  '

$ws.Range("D50").Value = '''25.69'
$ws.Range("E50").Value = '  +0.20%  '

$ws.Range("E51").Value = '  -0.09%  '
